# Auto update Excel log
# Appends new sensor readings to the "PIR" sheet (rows 81-93) and the
# "Humidity" sheet (rows 56-64).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# PIR sheet: append "No Motion" / "Inactive" bathroom readings
# ---------------------------------------------------------------------
$pir = $wb.Worksheets.Item("PIR")

$pirTimestamps = @(
    "18:25:16","18:25:18","18:25:23","18:25:28","18:25:33","18:25:38",
    "18:25:43","18:25:48","18:25:53","18:25:58","18:26:03","18:26:08","18:26:13"
)

$pirStartRow = 81
$pirEndRow = $pirStartRow + $pirTimestamps.Length - 1

# Keep the Date column as plain text (matches the rest of the log) instead
# of letting Excel auto-convert "2026-01-30" into a date serial number.
$pir.Range("A$pirStartRow`:A$pirEndRow").NumberFormat = "@"

for ($i = 0; $i -lt $pirTimestamps.Length; $i++) {
    $r = $pirStartRow + $i
    $pir.Cells.Item($r, 1).Value = "2026-01-30"
    $pir.Cells.Item($r, 2).Value = $pirTimestamps[$i]
    $pir.Cells.Item($r, 3).Value = "18:00"
    $pir.Cells.Item($r, 4).Value = "Bathroom"
    $pir.Cells.Item($r, 5).Value = "No Motion"
    $pir.Cells.Item($r, 6).Value = "Inactive"
}

# ---------------------------------------------------------------------
# Humidity sheet: append "Active" bathroom humidity readings
# ---------------------------------------------------------------------
$hum = $wb.Worksheets.Item("Humidity")

$humRows = @(
    @("18:25:17","85.9%"),
    @("18:25:23","86.8%"),
    @("18:25:28","86.7%"),
    @("18:25:33","85.8%"),
    @("18:25:43","86.7%"),
    @("18:25:48","86.7%"),
    @("18:25:53","85.8%"),
    @("18:26:08","86.6%"),
    @("18:26:13","85.7%")
)

$humStartRow = 56
$humEndRow = $humStartRow + $humRows.Length - 1

# Keep Date and Value (percentage) columns as plain text, matching the
# rest of the log, instead of letting Excel auto-convert them into a
# date serial number / percentage number.
$hum.Range("A$humStartRow`:A$humEndRow").NumberFormat = "@"
$hum.Range("E$humStartRow`:E$humEndRow").NumberFormat = "@"

for ($i = 0; $i -lt $humRows.Length; $i++) {
    $r = $humStartRow + $i
    $hum.Cells.Item($r, 1).Value = "2026-01-30"
    $hum.Cells.Item($r, 2).Value = $humRows[$i][0]
    $hum.Cells.Item($r, 3).Value = "18:00"
    $hum.Cells.Item($r, 4).Value = "Bathroom"
    $hum.Cells.Item($r, 5).Value = $humRows[$i][1]
    $hum.Cells.Item($r, 6).Value = "Active"
}
